$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price-report row was added to the top of this week's block
# (row 68), pushing every existing record from row 68 down to row 69, and
# so on through the former last row (124 -> 125).
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new record's data.
$ws.Cells.Item(68, 1).Value  = 11
$ws.Cells.Item(68, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(68, 3).Value  = "Bíobío"
$ws.Cells.Item(68, 4).Value  = 44484
$ws.Cells.Item(68, 5).Value  = 8
$ws.Cells.Item(68, 6).Value  = 100114001
$ws.Cells.Item(68, 7).Value  = "Papa"
$ws.Cells.Item(68, 8).Value  = "Asterix"
$ws.Cells.Item(68, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(68, 10).Value = 650
$ws.Cells.Item(68, 11).Value = 11000
$ws.Cells.Item(68, 12).Value = 12000
$ws.Cells.Item(68, 13).Value = 11462
$ws.Cells.Item(68, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(68, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(68, 16).Value = 458
$ws.Cells.Item(68, 17).Value = 25
$ws.Cells.Item(68, 18).Value = "Hortaliza"
